# Update Name of Algo
# Apply targeted numeric corrections to column B (and one in column C)
# as produced by a re-run of the RandomForest imputation algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value  = 4.7807
$ws.Range("B10").Value = 8.639600000000003
$ws.Range("B12").Value = 5.739900000000002
$ws.Range("C13").Value = -12.59679999999999
$ws.Range("B18").Value = 5.144200000000004
